$d = $word.ActiveDocument

# The document's single table gets one new row appended at the bottom,
# continuing the existing "dicionario de dados" numbering (item 14).
$table = $d.Tables.Item(1)

$newRow = $table.Rows.Add()
$rowIndex = $newRow.Index

$table.Cell($rowIndex, 1).Range.Text = "14."
$table.Cell($rowIndex, 2).Range.Text = "LINK_PORTAL"
$table.Cell($rowIndex, 3).Range.Text = "URL"
$table.Cell($rowIndex, 4).Range.Text = "Link para o portal da transparência para os processos de compras que não possuem contratos, que são processos com entrega imediata."
